$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.138.91'
$ws.Range("E2").Value = '  +5.77%  '
$ws.Range("D3").Value = '1.923.91'
$ws.Range("E3").Value = '  +2.70%  '
$ws.Range("E4").Value = '  -0.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.59'
$ws.Range("E5").Value = '  +5.03%  '
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5238'
$ws.Range("E7").Value = '  +2.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4109'
$ws.Range("E8").Value = '  +5.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08549'
$ws.Range("E9").Value = '  +2.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.43'
$ws.Range("E10").Value = '  +3.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.131'
$ws.Range("E11").Value = '  +2.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.48'
$ws.Range("E12").Value = '  +10.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.451'
$ws.Range("E13").Value = '  +3.64%  '
$ws.Range("D14").Value = '1.921.14'
$ws.Range("E14").Value = '  +2.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.424'
$ws.Range("E15").Value = '  +2.15%  '
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.60'
$ws.Range("E17").Value = '  +5.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001118'
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06717'
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.36'
$ws.Range("E20").Value = '  +3.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.070'
$ws.Range("E22").Value = '  +2.77%  '
$ws.Range("D23").Value = '30.158.60'
$ws.Range("E23").Value = '  +5.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.34'
$ws.Range("E24").Value = '  +1.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.217'
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("D26").Value = '2.141.48'
$ws.Range("E26").Value = '  +2.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.20'
$ws.Range("E27").Value = '  +2.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.93'
$ws.Range("E28").Value = '  -1.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.477'
$ws.Range("E29").Value = '  +1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.72'
$ws.Range("E30").Value = '  +2.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.084'
$ws.Range("E31").Value = '  +4.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1057'
$ws.Range("E32").Value = '  +1.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.130'
$ws.Range("E33").Value = '  +6.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.646'
$ws.Range("E34").Value = '  +1.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02518'
$ws.Range("E35").Value = '  +2.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06614'
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("E37").Value = '  +3.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.243'
$ws.Range("E38").Value = '  +4.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.044'
$ws.Range("E39").Value = '  +2.32%  '
$ws.Range("E40").Value = '  +3.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6566'
$ws.Range("E41").Value = '  +3.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.67'
$ws.Range("E42").Value = '  +5.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.244'
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6195'
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.30'
$ws.Range("E45").Value = '  +2.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.778'
$ws.Range("E46").Value = '  +2.38%  '
$ws.Range("E47").Value = '  +4.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.249'
$ws.Range("E48").Value = '  +2.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.28'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.162'
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.08'
$ws.Range("E51").Value = '  +5.14%  '
